$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Disposition" status of the latest action items ---

# Row 15: "Talk to the UPT team..." -> more detailed Complete note
$ws.Range("E15").Value = "Complete (UPT team to prioritize and schedule the 2-3 week effort.)"

# Rows 16 & 17: "Follow up with Doug Hosier..." / "Follow up on CIT Security..." -> Complete
$ws.Range("E16").Value = "Complete"
$ws.Range("E17").Value = "Complete"

# --- Normalize cell styling for rows 12, 16 and 17 to match the rest of the table ---
# (copy number format / font / fill / border from the plain row style used elsewhere)
$ws.Range("A13:E13").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A16:E17").PasteSpecial(-4122)

# Rows 16 & 17 grew a line of wrapped text worth of height, same as other wrapped rows
$ws.Rows(16).RowHeight = 31
$ws.Rows(17).RowHeight = 31

# Leave the freshly-updated rows selected, like the author would after editing them
$ws.Range("A16:XFD17").Select()

$excel.CutCopyMode = 0
